$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 343.5  # H12: 331 -> 343.5
$ws.Cells.Item(12, 9).Value = 349.5  # I12: 299.6 -> 349.5
$ws.Cells.Item(12, 10).Value = 337.5  # J12: 383.33334 -> 337.5
$ws.Cells.Item(12, 11).Value = 349.5  # K12: 299.6 -> 349.5
$ws.Cells.Item(12, 12).Value = 337.5  # L12: 383.33334 -> 337.5
$ws.Cells.Item(12, 13).Value = -179.5  # M12: -129.6 -> -179.5
$ws.Cells.Item(12, 14).Value = -677.5  # N12: -723.33334 -> -677.5

$ws.Cells.Item(21, 8).Value = 23952.334  # H21: 20000 -> 23952.334
$ws.Cells.Item(21, 9).Value = 0  # I21: 20000 -> 0
$ws.Cells.Item(21, 10).Value = 23952.334  # J21: 0 -> 23952.334
$ws.Cells.Item(21, 11).Value = 0  # K21: 20000 -> 0
$ws.Cells.Item(21, 12).Value = 23952.334  # L21: 0 -> 23952.334
$ws.Cells.Item(21, 13).ClearContents()  # M21: -19532 -> (removed)
$ws.Cells.Item(21, 14).Value = -24888.334  # N21: None -> -24888.334

$ws.Cells.Item(23, 8).Value = 23952.334  # H23: 20000 -> 23952.334
$ws.Cells.Item(23, 9).Value = 0  # I23: 20000 -> 0
$ws.Cells.Item(23, 10).Value = 23952.334  # J23: 0 -> 23952.334
$ws.Cells.Item(23, 11).Value = 0  # K23: 20000 -> 0
$ws.Cells.Item(23, 12).Value = 23952.334  # L23: 0 -> 23952.334
$ws.Cells.Item(23, 13).ClearContents()  # M23: -19766 -> (removed)
$ws.Cells.Item(23, 14).Value = -24420.334  # N23: None -> -24420.334

$ws.Cells.Item(29, 8).Value = 1700  # H29: 1795 -> 1700
$ws.Cells.Item(29, 10).Value = 2228.5715  # J29: 2320 -> 2228.5715
$ws.Cells.Item(29, 12).Value = 6685.7145  # L29: 6960 -> 6685.7145
$ws.Cells.Item(29, 14).Value = -7247.7145  # N29: -7522 -> -7247.7145

$ws.Cells.Item(38, 8).Value = 1625.5853  # H38: 1657.4286 -> 1625.5853
$ws.Cells.Item(38, 9).Value = 164.125  # I38: 184 -> 164.125
$ws.Cells.Item(38, 10).Value = 1979.8788  # J38: 1952.1143 -> 1979.8788
$ws.Cells.Item(38, 11).Value = 492.375  # K38: 552 -> 492.375
$ws.Cells.Item(38, 12).Value = 5939.636399999999  # L38: 5856.3429 -> 5939.636399999999
$ws.Cells.Item(38, 13).Value = -120.375  # M38: -180 -> -120.375
$ws.Cells.Item(38, 14).Value = -6683.636399999999  # N38: -6600.3429 -> -6683.636399999999

$ws.Cells.Item(58, 8).Value = 1273.1765  # H58: 1156.8889 -> 1273.1765
$ws.Cells.Item(58, 9).Value = 493.7  # I58: 417.25 -> 493.7
$ws.Cells.Item(58, 10).Value = 2386.7144  # J58: 2636.1667 -> 2386.7144
$ws.Cells.Item(58, 11).Value = 1481.1  # K58: 1251.75 -> 1481.1
$ws.Cells.Item(58, 12).Value = 7160.1432  # L58: 7908.500100000001 -> 7160.1432
$ws.Cells.Item(58, 13).Value = -1331.1  # M58: -1101.75 -> -1331.1
$ws.Cells.Item(58, 14).Value = -7460.1432  # N58: -8208.500100000001 -> -7460.1432

$ws.Cells.Item(92, 8).Value = 2746.9092  # H92: 2241.6155 -> 2746.9092
$ws.Cells.Item(92, 9).Value = 3423.2  # I92: 2652.6667 -> 3423.2
$ws.Cells.Item(92, 10).Value = 2183.3333  # J92: 1889.2858 -> 2183.3333
$ws.Cells.Item(92, 11).Value = 3423.2  # K92: 2652.6667 -> 3423.2
$ws.Cells.Item(92, 12).Value = 2183.3333  # L92: 1889.2858 -> 2183.3333
$ws.Cells.Item(92, 13).Value = -2175.2  # M92: -1404.6667 -> -2175.2
$ws.Cells.Item(92, 14).Value = -4679.3333  # N92: -4385.2858 -> -4679.3333

$ws.Cells.Item(116, 8).Value = 2735  # H116: 2693.25 -> 2735
$ws.Cells.Item(116, 9).Value = 2183.7693  # I116: 2163.5 -> 2183.7693
$ws.Cells.Item(116, 11).Value = 2183.7693  # K116: 2163.5 -> 2183.7693
$ws.Cells.Item(116, 13).Value = 1258.2307  # M116: 1278.5 -> 1258.2307

$ws.Cells.Item(138, 8).Value = 436434.44  # H138: 440838.16 -> 436434.44
$ws.Cells.Item(138, 9).Value = 915.9524  # I138: 936.95 -> 915.9524
$ws.Cells.Item(138, 10).Value = 552205.2  # J138: 552205.5600000001 -> 552205.2
$ws.Cells.Item(138, 11).Value = 2747.8572  # K138: 2810.85 -> 2747.8572
$ws.Cells.Item(138, 12).Value = 1656615.6  # L138: 1656616.68 -> 1656615.6
$ws.Cells.Item(138, 13).Value = 2392.1428  # M138: 2329.15 -> 2392.1428
$ws.Cells.Item(138, 14).Value = -1666895.6  # N138: -1666896.68 -> -1666895.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6128.2607  # H32: 6757.524 -> 6128.2607
$ws.Cells.Item(32, 9).Value = 6042.2446  # I32: 6678.4634 -> 6042.2446
$ws.Cells.Item(32, 11).Value = 6042.2446  # K32: 6678.4634 -> 6042.2446
$ws.Cells.Item(32, 13).Value = -5755.2446  # M32: -6391.4634 -> -5755.2446

$ws.Cells.Item(74, 8).Value = 1018.7  # H74: 1002.6667 -> 1018.7
$ws.Cells.Item(74, 9).Value = 1018.7  # I74: 957.4545000000001 -> 1018.7
$ws.Cells.Item(74, 10).Value = 0  # J74: 1500 -> 0
$ws.Cells.Item(74, 11).Value = 1018.7  # K74: 957.4545000000001 -> 1018.7
$ws.Cells.Item(74, 12).Value = 0  # L74: 1500 -> 0
$ws.Cells.Item(74, 13).Value = -144.7  # M74: -83.45450000000005 -> -144.7
$ws.Cells.Item(74, 14).ClearContents()  # N74: -3248 -> (removed)

$ws.Cells.Item(77, 8).Value = 1018.7  # H77: 1002.6667 -> 1018.7
$ws.Cells.Item(77, 9).Value = 1018.7  # I77: 957.4545000000001 -> 1018.7
$ws.Cells.Item(77, 10).Value = 0  # J77: 1500 -> 0
$ws.Cells.Item(77, 11).Value = 5093.5  # K77: 4787.2725 -> 5093.5
$ws.Cells.Item(77, 12).Value = 0  # L77: 7500 -> 0
$ws.Cells.Item(77, 13).Value = -725.5  # M77: -419.2725 -> -725.5
$ws.Cells.Item(77, 14).ClearContents()  # N77: -16236 -> (removed)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6031.7144  # H134: 7739.625 -> 6031.7144
$ws.Cells.Item(134, 9).Value = 870.3333  # I134: 988.1429000000001 -> 870.3333
$ws.Cells.Item(134, 10).Value = 37000  # J134: 55000 -> 37000
$ws.Cells.Item(134, 11).Value = 2610.9999  # K134: 2964.4287 -> 2610.9999
$ws.Cells.Item(134, 12).Value = 111000  # L134: 165000 -> 111000
$ws.Cells.Item(134, 13).Value = -75.9998999999998  # M134: -429.4287000000004 -> -75.9998999999998
$ws.Cells.Item(134, 14).Value = -116070  # N134: -170070 -> -116070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 307.8  # H7: 288.41666 -> 307.8
$ws.Cells.Item(7, 9).Value = 127.6  # I7: 114.5 -> 127.6
$ws.Cells.Item(7, 10).Value = 488  # J7: 462.33334 -> 488
$ws.Cells.Item(7, 11).Value = 127.6  # K7: 114.5 -> 127.6
$ws.Cells.Item(7, 12).Value = 488  # L7: 462.33334 -> 488
$ws.Cells.Item(7, 13).Value = -14.59999999999999  # M7: -1.5 -> -14.59999999999999
$ws.Cells.Item(7, 14).Value = -714  # N7: -688.33334 -> -714

$ws.Cells.Item(31, 8).Value = 1710.3793  # H31: 1602.7812 -> 1710.3793
$ws.Cells.Item(31, 9).Value = 2169.7273  # I31: 1765.9375 -> 2169.7273
$ws.Cells.Item(31, 10).Value = 1429.6666  # J31: 1439.625 -> 1429.6666
$ws.Cells.Item(31, 11).Value = 2169.7273  # K31: 1765.9375 -> 2169.7273
$ws.Cells.Item(31, 12).Value = 1429.6666  # L31: 1439.625 -> 1429.6666
$ws.Cells.Item(31, 13).Value = -1874.7273  # M31: -1470.9375 -> -1874.7273
$ws.Cells.Item(31, 14).Value = -2019.6666  # N31: -2029.625 -> -2019.6666

$ws.Cells.Item(34, 8).Value = 1710.3793  # H34: 1602.7812 -> 1710.3793
$ws.Cells.Item(34, 9).Value = 2169.7273  # I34: 1765.9375 -> 2169.7273
$ws.Cells.Item(34, 10).Value = 1429.6666  # J34: 1439.625 -> 1429.6666
$ws.Cells.Item(34, 11).Value = 2169.7273  # K34: 1765.9375 -> 2169.7273
$ws.Cells.Item(34, 12).Value = 1429.6666  # L34: 1439.625 -> 1429.6666
$ws.Cells.Item(34, 13).Value = -1967.7273  # M34: -1563.9375 -> -1967.7273
$ws.Cells.Item(34, 14).Value = -1833.6666  # N34: -1843.625 -> -1833.6666

$ws.Cells.Item(58, 8).Value = 1446.2222  # H58: 1539.7084 -> 1446.2222
$ws.Cells.Item(58, 9).Value = 1123.2222  # I58: 1208.2 -> 1123.2222
$ws.Cells.Item(58, 11).Value = 1123.2222  # K58: 1208.2 -> 1123.2222
$ws.Cells.Item(58, 13).Value = -920.2221999999999  # M58: -1005.2 -> -920.2221999999999

$ws.Cells.Item(132, 8).Value = 2739.5  # H132: 2377.3333 -> 2739.5
$ws.Cells.Item(132, 9).Value = 2614.4  # I132: 2133.5 -> 2614.4
$ws.Cells.Item(132, 10).Value = 2948  # J132: 3157.6 -> 2948
$ws.Cells.Item(132, 11).Value = 7843.200000000001  # K132: 6400.5 -> 7843.200000000001
$ws.Cells.Item(132, 12).Value = 8844  # L132: 9472.799999999999 -> 8844
$ws.Cells.Item(132, 13).Value = -5313.200000000001  # M132: -3870.5 -> -5313.200000000001
$ws.Cells.Item(132, 14).Value = -13904  # N132: -14532.8 -> -13904

$ws.Cells.Item(136, 8).Value = 1446.2222  # H136: 1539.7084 -> 1446.2222
$ws.Cells.Item(136, 9).Value = 1123.2222  # I136: 1208.2 -> 1123.2222
$ws.Cells.Item(136, 11).Value = 3369.6666  # K136: 3624.6 -> 3369.6666
$ws.Cells.Item(136, 13).Value = -819.6665999999996  # M136: -1074.6 -> -819.6665999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 3258.7058  # H39: 3400.125 -> 3258.7058
$ws.Cells.Item(39, 10).Value = 3186.5334  # J39: 3343 -> 3186.5334
$ws.Cells.Item(39, 12).Value = 9559.600199999999  # L39: 10029 -> 9559.600199999999
$ws.Cells.Item(39, 14).Value = -10147.6002  # N39: -10617 -> -10147.6002

$ws.Cells.Item(55, 8).Value = 1814.8572  # H55: 2042.3334 -> 1814.8572
$ws.Cells.Item(55, 10).Value = 2041.6666  # J55: 2360 -> 2041.6666
$ws.Cells.Item(55, 12).Value = 6124.9998  # L55: 7080 -> 6124.9998
$ws.Cells.Item(55, 14).Value = -6478.9998  # N55: -7434 -> -6478.9998

$ws.Cells.Item(131, 8).Value = 37039044  # H131: 30304798 -> 37039044
$ws.Cells.Item(131, 9).Value = 250000300  # I131: 200000320 -> 250000300
$ws.Cells.Item(131, 10).Value = 2305.1738  # J131: 2025.5 -> 2305.1738
$ws.Cells.Item(131, 11).Value = 750000900  # K131: 600000960 -> 750000900
$ws.Cells.Item(131, 12).Value = 6915.5214  # L131: 6076.5 -> 6915.5214
$ws.Cells.Item(131, 13).Value = -749995860  # M131: -599995920 -> -749995860
$ws.Cells.Item(131, 14).Value = -16995.5214  # N131: -16156.5 -> -16995.5214

$ws.Cells.Item(139, 8).Value = 1659.4324  # H139: 1627.079 -> 1659.4324
$ws.Cells.Item(139, 9).Value = 1635.5652  # I139: 1585.3334 -> 1635.5652
$ws.Cells.Item(139, 11).Value = 4906.6956  # K139: 4756.0002 -> 4906.6956
$ws.Cells.Item(139, 13).Value = 233.3044  # M139: 383.9997999999996 -> 233.3044

$ws.Cells.Item(140, 8).Value = 25165.777  # H140: 24222.426 -> 25165.777
$ws.Cells.Item(140, 10).Value = 3463.0356  # J140: 3431.9666 -> 3463.0356
$ws.Cells.Item(140, 12).Value = 10389.1068  # L140: 10295.8998 -> 10389.1068
$ws.Cells.Item(140, 14).Value = -20749.1068  # N140: -20655.8998 -> -20749.1068

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1883.4412  # H122: 2266.5806 -> 1883.4412
$ws.Cells.Item(122, 9).Value = 1838  # I122: 2348.4783 -> 1838
$ws.Cells.Item(122, 11).Value = 5514  # K122: 7045.4349 -> 5514
$ws.Cells.Item(122, 13).Value = -3064  # M122: -4595.4349 -> -3064

$ws.Cells.Item(132, 8).Value = 5139.273  # H132: 6377.75 -> 5139.273
$ws.Cells.Item(132, 9).Value = 6904.2  # I132: 14505.5 -> 6904.2
$ws.Cells.Item(132, 11).Value = 20712.6  # K132: 43516.5 -> 20712.6
$ws.Cells.Item(132, 13).Value = -18182.6  # M132: -40986.5 -> -18182.6

$ws.Cells.Item(136, 8).Value = 13208.3125  # H136: 11078.412 -> 13208.3125
$ws.Cells.Item(136, 10).Value = 13208.3125  # J136: 11078.412 -> 13208.3125
$ws.Cells.Item(136, 12).Value = 39624.9375  # L136: 33235.236 -> 39624.9375
$ws.Cells.Item(136, 14).Value = -44724.9375  # N136: -38335.236 -> -44724.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5425.7144  # H46: 4553.3335 -> 5425.7144
$ws.Cells.Item(46, 9).Value = 1660  # I46: 1745 -> 1660
$ws.Cells.Item(46, 10).Value = 8250  # J46: 6800 -> 8250
$ws.Cells.Item(46, 11).Value = 1660  # K46: 1745 -> 1660
$ws.Cells.Item(46, 12).Value = 8250  # L46: 6800 -> 8250
$ws.Cells.Item(46, 13).Value = -1472  # M46: -1557 -> -1472
$ws.Cells.Item(46, 14).Value = -8626  # N46: -7176 -> -8626

$ws.Cells.Item(132, 8).Value = 2874.1304  # H132: 2804.2 -> 2874.1304
$ws.Cells.Item(132, 9).Value = 2845.182  # I132: 2666.4167 -> 2845.182
$ws.Cells.Item(132, 10).Value = 2900.6667  # J132: 2931.3845 -> 2900.6667
$ws.Cells.Item(132, 11).Value = 8535.545999999998  # K132: 7999.250100000001 -> 8535.545999999998
$ws.Cells.Item(132, 12).Value = 8702.000100000001  # L132: 8794.1535 -> 8702.000100000001
$ws.Cells.Item(132, 13).Value = -6005.545999999998  # M132: -5469.250100000001 -> -6005.545999999998
$ws.Cells.Item(132, 14).Value = -13762.0001  # N132: -13854.1535 -> -13762.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 75011  # H20: 83344.336 -> 75011
$ws.Cells.Item(20, 10).Value = 75011  # J20: 83344.336 -> 75011
$ws.Cells.Item(20, 12).Value = 75011  # L20: 83344.336 -> 75011
$ws.Cells.Item(20, 14).Value = -75491  # N20: -83824.336 -> -75491

$ws.Cells.Item(62, 8).Value = 125006250  # H62: 83338340 -> 125006250
$ws.Cells.Item(62, 10).Value = 8334.333000000001  # J62: 6000.6 -> 8334.333000000001
$ws.Cells.Item(62, 12).Value = 8334.333000000001  # L62: 6000.6 -> 8334.333000000001
$ws.Cells.Item(62, 14).Value = -9582.333000000001  # N62: -7248.6 -> -9582.333000000001

$ws.Cells.Item(65, 8).Value = 125006250  # H65: 83338340 -> 125006250
$ws.Cells.Item(65, 10).Value = 8334.333000000001  # J65: 6000.6 -> 8334.333000000001
$ws.Cells.Item(65, 12).Value = 41671.665  # L65: 30003 -> 41671.665
$ws.Cells.Item(65, 14).Value = -47911.665  # N65: -36243 -> -47911.665

$ws.Cells.Item(81, 8).Value = 375.125  # H81: 740 -> 375.125
$ws.Cells.Item(81, 9).Value = 285.85715  # I81: 450 -> 285.85715
$ws.Cells.Item(81, 10).Value = 1000  # J81: 933.3333 -> 1000
$ws.Cells.Item(81, 11).Value = 571.7143  # K81: 900 -> 571.7143
$ws.Cells.Item(81, 12).Value = 2000  # L81: 1866.6666 -> 2000
$ws.Cells.Item(81, 13).Value = 489.2857  # M81: 161 -> 489.2857
$ws.Cells.Item(81, 14).Value = -4122  # N81: -3988.6666 -> -4122

$ws.Cells.Item(84, 8).Value = 375.125  # H84: 740 -> 375.125
$ws.Cells.Item(84, 9).Value = 285.85715  # I84: 450 -> 285.85715
$ws.Cells.Item(84, 10).Value = 1000  # J84: 933.3333 -> 1000
$ws.Cells.Item(84, 11).Value = 2858.5715  # K84: 4500 -> 2858.5715
$ws.Cells.Item(84, 12).Value = 10000  # L84: 9333.333000000001 -> 10000
$ws.Cells.Item(84, 13).Value = 2445.4285  # M84: 804 -> 2445.4285
$ws.Cells.Item(84, 14).Value = -20608  # N84: -19941.333 -> -20608

$ws.Cells.Item(122, 8).Value = 25001996  # H122: 15626489 -> 25001996
$ws.Cells.Item(122, 9).Value = 31252242  # I122: 20835044 -> 31252242
$ws.Cells.Item(122, 10).Value = 1005  # J122: 825 -> 1005
$ws.Cells.Item(122, 11).Value = 93756726  # K122: 62505132 -> 93756726
$ws.Cells.Item(122, 12).Value = 3015  # L122: 2475 -> 3015
$ws.Cells.Item(122, 13).Value = -93754276  # M122: -62502682 -> -93754276
$ws.Cells.Item(122, 14).Value = -7915  # N122: -7375 -> -7915
